$wb = $excel.ActiveWorkbook

# --- Update the conversion text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.58 = 51132.08 pesos`n✅ 51132.08 pesos = 12.51 = 972.34 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 79.5
$wsTasas.Range("O10").Value = 4065
$wsTasas.Range("N12").Value = 4085.99
$wsTasas.Range("O12").Value = 77.7
